# Hortaliza, Vega Monumental Concepción - Brócoli
# Insert a new weekly price-report row at row 74 (pushing the existing
# rows 74-159 down to 75-160), matching the "Fruta / hortaliza, semanal"
# commit that adds one more sample for Brócoli.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 74; everything below shifts down by one.
$ws.Rows("74").Insert()

# Populate the new row 74 with the new weekly record.
$ws.Cells.Item(74, 1).Value = 11
$ws.Cells.Item(74, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(74, 3).Value = "Bíobío"
$ws.Cells.Item(74, 4).Value = 44483
$ws.Cells.Item(74, 5).Value = 8
$ws.Cells.Item(74, 6).Value = 100112023
$ws.Cells.Item(74, 7).Value = "Brócoli"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 2100
$ws.Cells.Item(74, 11).Value = 600
$ws.Cells.Item(74, 12).Value = 650
$ws.Cells.Item(74, 13).Value = 629
$ws.Cells.Item(74, 14).Value = "$/unidad"
$ws.Cells.Item(74, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(74, 16).Value = 629
$ws.Cells.Item(74, 17).Value = 1
$ws.Cells.Item(74, 18).Value = "Hortaliza"
